$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing autofilter before the row shift so re-adding it below
# installs a clean record rather than toggling the stale one off.
$ws.AutoFilterMode = $false

# Delete the first row (descriptive functional-group header row, e.g.
# "Carbohydrate: ester bonds", etc.); everything below shifts up by one,
# turning the old row 2 (id/veg/prec/time/carbo1/...) into row 1.
$ws.Rows.Item(1).Delete()

# Re-apply the autofilter over the new header/data extent (was D2:D66,
# now D1:D65 after the shift).
$ws.Range("D1:D65").AutoFilter()

# The autofilter's backing defined name needs the same adjustment.
foreach ($n in $wb.Names) {
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=ftir.percent.area_2!`$D`$1:`$D`$65"
  }
}

# Reflect the new active selection as captured in the saved file.
$ws.Range("F19").Select()
